# "develope to date function" - trims the placeholder "X" marks that used to
# fill every weekday/holiday/consecutive-holiday cell (columns E/F/G) down to
# only the cells that are still meaningful, in preparation for a date-driven
# formula taking over those columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-10: drop the "平日" (E) marker, keep F/G
$rowsDropE = 5..10
foreach ($r in $rowsDropE) {
    $ws.Range("E$r").ClearContents()
}

# Rows 11-16: drop the "平日"/"假日" (E/F) markers, keep G
$rowsDropEF = 11..16
foreach ($r in $rowsDropEF) {
    $ws.Range("E$r").ClearContents()
    $ws.Range("F$r").ClearContents()
}

# Rows 17-29: drop the "平日"/"連續假日" (E/G) markers entirely (F already empty)
$rowsDropEG = 17..29
foreach ($r in $rowsDropEG) {
    $ws.Range("E$r").ClearContents()
    $ws.Range("G$r").ClearContents()
}

# Rows 30-32: drop the remaining "連續假日" (G) marker
$rowsDropG = 30..32
foreach ($r in $rowsDropG) {
    $ws.Range("G$r").ClearContents()
}

# Move the active selection from L7 to I10
$ws.Range("I10").Select()
